$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set text number format on Price (D) and Volume(1h) (E) columns for the edited cells
# so that numeric-looking strings (e.g. "1.005") are preserved as text, matching
# the source data which stores these as plain strings, not numbers.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '30.457.86'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +0.20%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.107.91'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +0.23%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.005'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.80%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '334.79'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +1.72%  '
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +0.73%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5219'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -0.30%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.4536'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +4.21%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '54.35'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +16.14%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.08924'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +0.85%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.182'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +1.51%  '
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -1.93%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '2.103.81'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +0.72%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.818'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +1.21%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '8.023'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +3.25%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '96.77'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +0.25%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001146'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +1.36%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '1.005'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +0.55%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06645'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +0.20%  '
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +1.14%  '
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +0.33%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.316'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -0.65%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '30.524.71'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +0.28%  '
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +0.21%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.347'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +0.99%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.349.84'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +0.68%  '
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -1.26%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '162.79'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +0.73%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.533'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -2.90%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '133.81'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +1.12%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.209'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -0.14%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.1068'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -0.42%  '
$ws.Range('B33').Value = 'ARBITRUM'
$ws.Range('C33').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.646'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -4.26%  '
$ws.Range('B34').Value = 'Filecoin'
$ws.Range('C34').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '6.385'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +3.06%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.942'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +1.23%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '10.37'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +3.41%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.776'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +5.21%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.02582'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -0.26%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.06837'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +1.90%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.2305'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +1.71%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '12.74'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +0.37%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.6872'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +0.62%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.247'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -0.53%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.321'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +4.83%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '13.96'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -0.59%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.6364'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -0.30%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.662'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.00000000353'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +23.38%  '
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -0.15%  '
$ws.Range('B50').Value = 'Aave'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '83.11'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +1.27%  '
$ws.Range('B51').Value = 'WEMIXTOKEN'
$ws.Range('C51').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.202'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +0.47%  '
